$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Select()
